# Applies the "added team spec time commit pt2" edit: updates probability/
# proportion values in the team-specific matrix (Sheet1) to match the new
# commit data, cell by cell, leaving row/column structure untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.194006309148265
$ws.Range("C2").Value = 0.5425867507886435
$ws.Range("J2").Value = 0.01261829652996845
$ws.Range("P2").Value = 0.1482649842271293
$ws.Range("S2").Value = 0.1025236593059937
$ws.Range("B3").Value = 0.005797101449275362
$ws.Range("C3").Value = 0.002898550724637681
$ws.Range("J3").Value = 0.01159420289855072
$ws.Range("P3").Value = 0.7304347826086957
$ws.Range("S3").Value = 0.2492753623188406
$ws.Range("J4").Value = 0.0505050505050505
$ws.Range("P4").Value = 0.696969696969697
$ws.Range("S4").Value = 0.2525252525252525
$ws.Range("B6").Value = 0.06720977596741344
$ws.Range("D6").Value = 0.01832993890020367
$ws.Range("F6").Value = 0.07535641547861507
$ws.Range("J6").Value = 0.2892057026476578
$ws.Range("O6").Value = 0.02240325865580448
$ws.Range("Q6").Value = 0.1384928716904277
$ws.Range("R6").Value = 0.06313645621181263
$ws.Range("S6").Value = 0.3258655804480652
$ws.Range("B7").Value = 0.1271186440677966
$ws.Range("D7").Value = 0.02542372881355932
$ws.Range("F7").Value = 0.03389830508474576
$ws.Range("J7").Value = 0.1299435028248588
$ws.Range("O7").Value = 0.03107344632768362
$ws.Range("Q7").Value = 0.1581920903954802
$ws.Range("R7").Value = 0.07344632768361582
$ws.Range("S7").Value = 0.4209039548022599
$ws.Range("B8").Value = 0.08986415882967608
$ws.Range("D8").Value = 0.02716823406478579
$ws.Range("E8").Value = 0.001044932079414838
$ws.Range("F8").Value = 0.06374085684430512
$ws.Range("J8").Value = 0.1024033437826541
$ws.Range("O8").Value = 0.01776384535005225
$ws.Range("Q8").Value = 0.1630094043887147
$ws.Range("R8").Value = 0.1086729362591431
$ws.Range("S8").Value = 0.4263322884012539
$ws.Range("B9").Value = 0.08123791102514506
$ws.Range("D9").Value = 0.01740812379110251
$ws.Range("F9").Value = 0.05029013539651837
$ws.Range("J9").Value = 0.1141199226305609
$ws.Range("O9").Value = 0.02707930367504836
$ws.Range("Q9").Value = 0.1972920696324952
$ws.Range("R9").Value = 0.1005802707930367
$ws.Range("S9").Value = 0.4119922630560928
$ws.Range("B10").Value = 0.1091820987654321
$ws.Range("D10").Value = 0.02006172839506173
$ws.Range("E10").Value = 0.0003858024691358024
$ws.Range("F10").Value = 0.06404320987654322
$ws.Range("J10").Value = 0.1226851851851852
$ws.Range("O10").Value = 0.01427469135802469
$ws.Range("Q10").Value = 0.1809413580246914
$ws.Range("R10").Value = 0.09027777777777778
$ws.Range("S10").Value = 0.3981481481481481
$ws.Range("G11").Value = 0.1474820143884892
$ws.Range("J11").Value = 0.1133093525179856
$ws.Range("K11").Value = 0.2050359712230216
$ws.Range("L11").Value = 0.512589928057554
$ws.Range("S11").Value = 0.02158273381294964
$ws.Range("G12").Value = 0.7389830508474576
$ws.Range("J12").Value = 0.1796610169491525
$ws.Range("K12").Value = 0.0135593220338983
$ws.Range("L12").Value = 0.03050847457627119
$ws.Range("S12").Value = 0.03728813559322034
$ws.Range("F13").Value = 0.0108695652173913
$ws.Range("G13").Value = 0.6413043478260869
$ws.Range("J13").Value = 0.2717391304347826
$ws.Range("S13").Value = 0.07608695652173914
$ws.Range("F15").Value = 0.0303030303030303
$ws.Range("H15").Value = 0.1535353535353535
$ws.Range("I15").Value = 0.07474747474747474
$ws.Range("J15").Value = 0.3515151515151515
$ws.Range("K15").Value = 0.04646464646464647
$ws.Range("M15").Value = 0.01414141414141414
$ws.Range("O15").Value = 0.07474747474747474
$ws.Range("S15").Value = 0.2545454545454545
$ws.Range("F16").Value = 0.02743142144638404
$ws.Range("H16").Value = 0.1745635910224439
$ws.Range("I16").Value = 0.08478802992518704
$ws.Range("J16").Value = 0.3865336658354115
$ws.Range("K16").Value = 0.09975062344139651
$ws.Range("M16").Value = 0.02244389027431421
$ws.Range("N16").Value = 0.004987531172069825
$ws.Range("O16").Value = 0.06234413965087282
$ws.Range("S16").Value = 0.1371571072319202
$ws.Range("F17").Value = 0.01995305164319249
$ws.Range("H17").Value = 0.1948356807511737
$ws.Range("I17").Value = 0.1009389671361502
$ws.Range("J17").Value = 0.3955399061032864
$ws.Range("K17").Value = 0.08568075117370892
$ws.Range("M17").Value = 0.01643192488262911
$ws.Range("N17").Value = 0.001173708920187793
$ws.Range("O17").Value = 0.06455399061032864
$ws.Range("S17").Value = 0.1208920187793427
$ws.Range("F18").Value = 0.0449438202247191
$ws.Range("H18").Value = 0.1685393258426966
$ws.Range("I18").Value = 0.1280898876404494
$ws.Range("J18").Value = 0.3370786516853932
$ws.Range("K18").Value = 0.09213483146067415
$ws.Range("M18").Value = 0.008988764044943821
$ws.Range("N18").Value = 0.002247191011235955
$ws.Range("O18").Value = 0.0898876404494382
$ws.Range("S18").Value = 0.1280898876404494
$ws.Range("F19").Value = 0.02479338842975207
$ws.Range("H19").Value = 0.2091268415379087
$ws.Range("I19").Value = 0.1092346388789077
$ws.Range("J19").Value = 0.3568091987064319
$ws.Range("K19").Value = 0.09019044196909809
$ws.Range("M19").Value = 0.02048149478979518
$ws.Range("N19").Value = 0.0007186489399928135
$ws.Range("O19").Value = 0.06970894717930291
$ws.Range("S19").Value = 0.1189363995688106
